# Update res_bus vm_pu results for the 380 kV slack-voltage case
# (slack setpoint B column 1.05 -> 1.02, downstream bus voltages refreshed)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.02897389898027
$ws.Range("D2").Value = 1.034483830102144
$ws.Range("E2").Value = 1.028851601151779
$ws.Range("F2").Value = 1.039863569151429
$ws.Range("I2").Value = 1.03715600899761
$ws.Range("J2").Value = 1.034123246083551
$ws.Range("K2").Value = 1.03728317289389
$ws.Range("L2").Value = 1.031667215805606
$ws.Range("M2").Value = 1.042647548265821
$ws.Range("N2").Value = 1.015312658946308

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.030109619414863
$ws.Range("D3").Value = 1.035071185316288
$ws.Range("E3").Value = 1.029821787742673
$ws.Range("F3").Value = 1.041198361270484
$ws.Range("I3").Value = 1.037450762013082
$ws.Range("J3").Value = 1.034898849031676
$ws.Range("K3").Value = 1.037680316524502
$ws.Range("L3").Value = 1.032444987654023
$ws.Range("M3").Value = 1.04379126516131
$ws.Range("N3").Value = 1.015575590033475

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030844168367823
$ws.Range("D4").Value = 1.035451031780925
$ws.Range("E4").Value = 1.030449608331414
$ws.Range("F4").Value = 1.042061959795502
$ws.Range("I4").Value = 1.03763994703523
$ws.Range("J4").Value = 1.0353998633757
$ws.Range("K4").Value = 1.037936385556385
$ws.Range("L4").Value = 1.032947690914953
$ws.Range("M4").Value = 1.04453069046437
$ws.Range("N4").Value = 1.015745284442494

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.031152893413881
$ws.Range("D5").Value = 1.035610667839212
$ws.Range("E5").Value = 1.030713555536255
$ws.Range("F5").Value = 1.042424994857565
$ws.Range("I5").Value = 1.037719111911104
$ws.Range("J5").Value = 1.035610286530234
$ws.Range("K5").Value = 1.038043818991152
$ws.Range("L5").Value = 1.033158892289576
$ws.Range("M5").Value = 1.044841394666026
$ws.Range("N5").Value = 1.015816518916882

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.031204725102187
$ws.Range("D6").Value = 1.035637468400601
$ws.Range("E6").Value = 1.030757874122669
$ws.Range("F6").Value = 1.042485948837975
$ws.Range("I6").Value = 1.037732382436814
$ws.Range("J6").Value = 1.035645605608092
$ws.Range("K6").Value = 1.038061844747067
$ws.Range("L6").Value = 1.033194346028951
$ws.Range("M6").Value = 1.044893554521074
$ws.Range("N6").Value = 1.015828473343374

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030848293876712
$ws.Range("D7").Value = 1.035453165048634
$ws.Range("E7").Value = 1.030453135161157
$ws.Range("F7").Value = 1.042066810770172
$ws.Range("I7").Value = 1.037641006287984
$ws.Range("J7").Value = 1.035402675857876
$ws.Range("K7").Value = 1.037937821944076
$ws.Range("L7").Value = 1.032950513527952
$ws.Range("M7").Value = 1.044534842696715
$ws.Range("N7").Value = 1.015746236693596

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029357791997296
$ws.Range("D8").Value = 1.034682372930811
$ws.Range("E8").Value = 1.029179471335095
$ws.Range("F8").Value = 1.04031469084993
$ws.Range("I8").Value = 1.037255941091551
$ws.Range("J8").Value = 1.034385541782976
$ws.Range("K8").Value = 1.037417577596221
$ws.Range("L8").Value = 1.031930185292506
$ws.Range("M8").Value = 1.043034205679076
$ws.Range("N8").Value = 1.015401608867983

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.026728675158388
$ws.Range("D9").Value = 1.033322549476359
$ws.Range("E9").Value = 1.026935427366159
$ws.Range("F9").Value = 1.037226345128282
$ws.Range("I9").Value = 1.036565606372323
$ws.Range("J9").Value = 1.03258665312403
$ws.Range("K9").Value = 1.036493886251125
$ws.Range("L9").Value = 1.030127861279484
$ws.Range("M9").Value = 1.040384927593377
$ws.Range("N9").Value = 1.014790955480023

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024974036753266
$ws.Range("D10").Value = 1.032414972792378
$ws.Range("E10").Value = 1.025439564438779
$ws.Range("F10").Value = 1.035166693979699
$ws.Range("I10").Value = 1.036097437732883
$ws.Range("J10").Value = 1.031382921715119
$ws.Range("K10").Value = 1.035873429337752
$ws.Range("L10").Value = 1.028923323441354
$ws.Range("M10").Value = 1.038615272502682
$ws.Range("N10").Value = 1.014381571193819

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.024213785812224
$ws.Range("D11").Value = 1.032021746347624
$ws.Range("E11").Value = 1.024791867832074
$ws.Range("F11").Value = 1.034274627567653
$ws.Range("I11").Value = 1.035892827893907
$ws.Range("J11").Value = 1.030860618729799
$ws.Range("K11").Value = 1.03560366087822
$ws.Range("L11").Value = 1.028401025256901
$ws.Range("M11").Value = 1.037848139173017
$ws.Range("N11").Value = 1.014203758849995

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.023931320342681
$ws.Range("D12").Value = 1.031875649103376
$ws.Range("E12").Value = 1.024551286725158
$ws.Range("F12").Value = 1.033943237559334
$ws.Range("I12").Value = 1.035816542437038
$ws.Range("J12").Value = 1.030666448750933
$ws.Range("K12").Value = 1.03550329071866
$ws.Range("L12").Value = 1.028206910424781
$ws.Range("M12").Value = 1.037563059642641
$ws.Range("N12").Value = 1.014137629048524

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.02399191354725
$ws.Range("D13").Value = 1.031906989052382
$ws.Range("E13").Value = 1.024602892104167
$ws.Range("F13").Value = 1.034014323558716
$ws.Range("I13").Value = 1.035832918795048
$ws.Range("J13").Value = 1.030708106284386
$ws.Range("K13").Value = 1.035524827975347
$ws.Range("L13").Value = 1.028248553711359
$ws.Range("M13").Value = 1.037624216188857
$ws.Range("N13").Value = 1.014151817848332

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.024190438659087
$ws.Range("D14").Value = 1.032009670623336
$ws.Range("E14").Value = 1.024771981282841
$ws.Range("F14").Value = 1.034247235535861
$ws.Range("I14").Value = 1.035886527914899
$ws.Range("J14").Value = 1.030844571919145
$ws.Range("K14").Value = 1.035595367638215
$ws.Range("L14").Value = 1.02838498191505
$ws.Range("M14").Value = 1.037824577119639
$ws.Range("N14").Value = 1.014198294222318

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024312746606405
$ws.Range("D15").Value = 1.032072931450086
$ws.Range("E15").Value = 1.024876162973565
$ws.Range("F15").Value = 1.034390735316983
$ws.Range("I15").Value = 1.035919520586464
$ws.Range("J15").Value = 1.030928631212211
$ws.Range("K15").Value = 1.035638807440036
$ws.Range("L15").Value = 1.028469025222184
$ws.Range("M15").Value = 1.037948008522331
$ws.Range("N15").Value = 1.014226918924718

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.02502448178705
$ws.Range("D16").Value = 1.032441064935389
$ws.Range("E16").Value = 1.025482550302446
$ws.Range("F16").Value = 1.035225892500784
$ws.Range("I16").Value = 1.03611097717051
$ws.Range("J16").Value = 1.031417562398279
$ws.Range("K16").Value = 1.035891309677504
$ws.Range("L16").Value = 1.02895797127946
$ws.Range("M16").Value = 1.038666166244802
$ws.Range("N16").Value = 1.014393360478067

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025470804177143
$ws.Range("D17").Value = 1.032671921573227
$ws.Range("E17").Value = 1.025862926188018
$ws.Range("F17").Value = 1.03574970235314
$ws.Range("I17").Value = 1.036230566653391
$ws.Range("J17").Value = 1.031723966173679
$ws.Range("K17").Value = 1.03604940148565
$ws.Range("L17").Value = 1.029264479283668
$ws.Range("M17").Value = 1.039116415194807
$ws.Range("N17").Value = 1.014497618401817

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.025731090218726
$ws.Range("D18").Value = 1.032806553173032
$ws.Range("E18").Value = 1.026084795191055
$ws.Range("F18").Value = 1.036055210391854
$ws.Range("I18").Value = 1.03630013886721
$ws.Range("J18").Value = 1.031902582254357
$ws.Range("K18").Value = 1.036141506992996
$ws.Range("L18").Value = 1.029443190302036
$ws.Range("M18").Value = 1.039378954897585
$ws.Range("N18").Value = 1.014558377601097

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.025819833252203
$ws.Range("D19").Value = 1.032852455121156
$ws.Range("E19").Value = 1.026160447189531
$ws.Range("F19").Value = 1.036159377207437
$ws.Range("I19").Value = 1.03632383028192
$ws.Range("J19").Value = 1.031963468109531
$ws.Range("K19").Value = 1.036172894477653
$ws.Range("L19").Value = 1.02950411431423
$ws.Range("M19").Value = 1.039468460173643
$ws.Range("N19").Value = 1.014579085991284

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.025422922808598
$ws.Range("D20").Value = 1.032647155233183
$ws.Range("E20").Value = 1.025822115236698
$ws.Range("F20").Value = 1.035693504747965
$ws.Range("I20").Value = 1.036217754693029
$ws.Range("J20").Value = 1.031691102720232
$ws.Range("K20").Value = 1.036032450773674
$ws.Range("L20").Value = 1.02923160111462
$ws.Range("M20").Value = 1.039068116336661
$ws.Range("N20").Value = 1.0144864379648

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024131980033536
$ws.Range("D21").Value = 1.031979434433466
$ws.Range("E21").Value = 1.0247221886804
$ws.Range("F21").Value = 1.034178649805947
$ws.Range("I21").Value = 1.0358707492148
$ws.Range("J21").Value = 1.030804390711472
$ws.Range("K21").Value = 1.035574600051223
$ws.Range("L21").Value = 1.028344810255349
$ws.Range("M21").Value = 1.037765579486487
$ws.Range("N21").Value = 1.014184610365183

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.023319881425939
$ws.Range("D22").Value = 1.031559407135922
$ws.Range("E22").Value = 1.024030633704437
$ws.Range("F22").Value = 1.033225985742812
$ws.Range("I22").Value = 1.035650928329072
$ws.Range("J22").Value = 1.030245933710849
$ws.Range("K22").Value = 1.035285769927054
$ws.Range("L22").Value = 1.027786612512488
$ws.Range("M22").Value = 1.036945857497495
$ws.Range("N22").Value = 1.013994362592594

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.02375043190355
$ws.Range("D23").Value = 1.03178209076657
$ws.Range("E23").Value = 1.024397239246001
$ws.Range("F23").Value = 1.033731032400765
$ws.Range("I23").Value = 1.035767615560137
$ws.Range("J23").Value = 1.030542072414522
$ws.Range("K23").Value = 1.035438975319969
$ws.Range("L23").Value = 1.02808258444446
$ws.Range("M23").Value = 1.037380481086424
$ws.Range("N23").Value = 1.014095261834897

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025444558480364
$ws.Range("D24").Value = 1.032658346147779
$ws.Range("E24").Value = 1.025840555941821
$ws.Range("F24").Value = 1.035718898092019
$ws.Range("I24").Value = 1.036223544429583
$ws.Range("J24").Value = 1.031705952621399
$ws.Range("K24").Value = 1.036040110400603
$ws.Range("L24").Value = 1.029246457559679
$ws.Range("M24").Value = 1.039089940767075
$ws.Range("N24").Value = 1.014491490085706

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.027408691038696
$ws.Range("D25").Value = 1.033674280527333
$ws.Range("E25").Value = 1.027515533869259
$ws.Range("F25").Value = 1.038024877817695
$ws.Range("I25").Value = 1.036745474040773
$ws.Range("J25").Value = 1.03305249310478
$ws.Range("K25").Value = 1.036733505362862
$ws.Range("L25").Value = 1.030594328528809
$ws.Range("M25").Value = 1.014949225294658
